# Updates cryptos list prices/volumes (and a few re-ranked rows) per the
# Mon Mar  6 02:36:02 UTC 2023 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '22.392.76'
$ws.Range("E2").Value = '  -0.63%  '

$ws.Range("D3").Value = '1.561.60'
$ws.Range("E3").Value = '  -0.85%  '

$ws.Range("D4").Value = '''1.002'
$ws.Range("E4").Value = '  +0.16%  '

$ws.Range("D5").Value = '''1.003'
$ws.Range("E5").Value = '  +0.18%  '

$ws.Range("D6").Value = '''286.13'
$ws.Range("E6").Value = '  -1.86%  '

$ws.Range("D7").Value = '''0.3622'
$ws.Range("E7").Value = '  -3.45%  '

$ws.Range("D8").Value = '''49.76'
$ws.Range("E8").Value = '  +0.42%  '

$ws.Range("D9").Value = '''0.3324'
$ws.Range("E9").Value = '  -2.37%  '

$ws.Range("D10").Value = '''1.120'
$ws.Range("E10").Value = '  -2.31%  '

$ws.Range("D11").Value = '''0.07365'
$ws.Range("E11").Value = '  -2.54%  '

$ws.Range("D12").Value = '''1.003'
$ws.Range("E12").Value = '  +0.27%  '

$ws.Range("D13").Value = '''20.79'
$ws.Range("E13").Value = '  -2.37%  '

$ws.Range("D14").Value = '''5.897'
$ws.Range("E14").Value = '  -1.62%  '

$ws.Range("D15").Value = '''6.851'
$ws.Range("E15").Value = '  -1.51%  '

$ws.Range("D16").Value = '1.566.89'
$ws.Range("E16").Value = '  -0.94%  '

$ws.Range("D17").Value = '''0.00001094'
$ws.Range("E17").Value = '  -2.37%  '

$ws.Range("D18").Value = '''88.83'
$ws.Range("E18").Value = '  -2.54%  '

$ws.Range("D19").Value = '''0.06722'
$ws.Range("E19").Value = '  -0.21%  '

$ws.Range("D20").Value = '''1.002'
$ws.Range("E20").Value = '  +0.20%  '

$ws.Range("D21").Value = '''6.274'
$ws.Range("E21").Value = '  +0.92%  '

$ws.Range("D22").Value = '''15.89'
$ws.Range("E22").Value = '  -3.10%  '

$ws.Range("D23").Value = '''11.89'
$ws.Range("E23").Value = '  -1.92%  '

$ws.Range("D24").Value = '22.384.79'
$ws.Range("E24").Value = '  -0.73%  '

$ws.Range("D25").Value = '''2.382'
$ws.Range("E25").Value = '  +0.39%  '

$ws.Range("D26").Value = '''2.502'
$ws.Range("E26").Value = '  -4.58%  '

$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").Value = '''148.78'
$ws.Range("E27").Value = '  +0.44%  '

$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").Value = '''19.54'
$ws.Range("E28").Value = '  -2.92%  '

$ws.Range("D29").Value = '''4.977'
$ws.Range("E29").Value = '  -0.05%  '

$ws.Range("D30").Value = '''122.71'
$ws.Range("E30").Value = '  -2.37%  '

$ws.Range("D31").Value = '1.737.53'
$ws.Range("E31").Value = '  -1.21%  '

$ws.Range("D32").Value = '''1.052'
$ws.Range("E32").Value = '  +0.95%  '

$ws.Range("D33").Value = '''6.065'
$ws.Range("E33").Value = '  -0.58%  '

$ws.Range("D34").Value = '''1.977'
$ws.Range("E34").Value = '  -0.64%  '

$ws.Range("D35").Value = '''9.476'
$ws.Range("E35").Value = '  -4.28%  '

$ws.Range("D36").Value = '''0.08243'
$ws.Range("E36").Value = '  -2.46%  '

$ws.Range("B37").Value = 'TrustWalletToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D37").Value = '''1.303'
$ws.Range("E37").Value = '  -5.58%  '

$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = '''0.02364'
$ws.Range("E38").Value = '  -4.18%  '

$ws.Range("D39").Value = '''0.2199'
$ws.Range("E39").Value = '  -3.42%  '

$ws.Range("D40").Value = '''0.06317'
$ws.Range("E40").Value = '  -3.66%  '

$ws.Range("D41").Value = '''5.266'
$ws.Range("E41").Value = '  -4.01%  '

$ws.Range("D42").Value = '''11.20'
$ws.Range("E42").Value = '  -1.86%  '

$ws.Range("B43").Value = 'Frax'
$ws.Range("C43").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D43").Value = '''1.001'
$ws.Range("E43").Value = '  +0.22%  '

$ws.Range("B44").Value = 'TheSandbox'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D44").Value = '''0.6018'
$ws.Range("E44").Value = '  -4.62%  '

$ws.Range("D45").Value = '''13.60'
$ws.Range("E45").Value = '  -4.14%  '

$ws.Range("D46").Value = '''3.743'
$ws.Range("E46").Value = '  -2.17%  '

$ws.Range("D47").Value = '''0.5665'
$ws.Range("E47").Value = '  -3.68%  '

$ws.Range("B48").Value = 'EOS'
$ws.Range("C48").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D48").Value = '''1.236'
$ws.Range("E48").Value = '  +0.58%  '

$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").Value = '''1.994'
$ws.Range("E49").Value = '  -4.96%  '

$ws.Range("B50").Value = 'Quant'
$ws.Range("C50").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D50").Value = '''123.65'
$ws.Range("E50").Value = '  -5.63%  '

$ws.Range("D51").Value = '''0.07248'
$ws.Range("E51").Value = '  -1.26%  '
